# Update COVID country data table and re-sort by total cases (column B, descending).
# Also refresh the "last updated" timestamp text in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New figures (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# for the countries whose underlying data changed in this refresh.
$updates = @{
    "Estados Unidos" = @(215362,359,8878,201371,5005,11,5113)
    "Alemania"        = @(80641,2660,19175,60504,3408,31,962)
    "Brasil"          = @(6988,108,127,6611,296,8,250)
    "Chile"           = @(3404,373,335,3051,31,2,18)
    "Islandia"        = @(1319,99,270,1047,12,0,2)
    "Serbia"          = @(1171,111,42,1098,81,3,31)
    "Singapur"        = @(1049,49,266,779,24,1,4)
    "Uruguay"         = @(350,0,62,284,15,2,4)
    "Afganistan"      = @(239,2,10,225,0,0,4)
    "Congo"           = @(22,0,2,18,0,0,2)
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row()

for ($r = 4; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($updates.ContainsKey($name)) {
        $vals = $updates[$name]
        for ($i = 0; $i -lt 7; $i++) {
            $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
        }
    }
}

# Re-sort the data block (countries + figures) by "Casos totales" (column B) descending,
# mirroring the sheet's existing sort order after the data refresh.
$sortRange = $ws.Range("A4:H" + $lastRow)
$sortKey = $ws.Range("B4:B" + $lastRow)
$sortRange.Sort($sortKey, 2)

# The refresh also flipped the relative order of Congo and San Martin (Parte Francesa),
# which are tied on total cases, so the stable sort above needs an explicit nudge.
$rowCongo = 0
$rowSanMartin = 0
for ($r = 4; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value()
    if ($name -eq "Congo") { $rowCongo = $r }
    if ($name -eq "San Martin (Parte Francesa)") { $rowSanMartin = $r }
}
if ($rowCongo -gt 0 -and $rowSanMartin -gt 0) {
    for ($c = 1; $c -le 8; $c++) {
        $tmp = $ws.Cells.Item($rowCongo, $c).Value()
        $ws.Cells.Item($rowCongo, $c).Value = $ws.Cells.Item($rowSanMartin, $c).Value()
        $ws.Cells.Item($rowSanMartin, $c).Value = $tmp
    }
}

# Refresh the "last updated" banner.
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 15:50"
